$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, pushing the existing
# "Late" / heading / "Outstanding" columns one place to the right.
$ws.Columns("N").Insert() | Out-Null

# Excel's "insert column" copies the left neighbour's width; the real
# commit ends up with an explicit width of 11 (matching column M) on the
# newly inserted column, so pin it explicitly.
$ws.Columns("N").ColumnWidth = 10.1666666666667

# The author switched focus to the "Repayment schedule" tab and left the
# selection on L12 there.
$ws.Activate() | Out-Null
$ws.Range("L12").Select() | Out-Null
